# Auto-generated Excel COM-interop script
# Applies updated market-price / leve-profit figures scraped by the scheduled runner
# to each job sheet's data table (columns H-N).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2196.9343
$ws.Range("I15").Value = 2196.9343
$ws.Range("K15").Value = 6590.8029
$ws.Range("M15").Value = -6421.8029
$ws.Range("H33").Value = 3313.5
$ws.Range("I33").Value = 3052.8572
$ws.Range("K33").Value = 3052.8572
$ws.Range("M33").Value = -2823.8572
$ws.Range("H47").Value = 25000
$ws.Range("I47").Value = 25000
$ws.Range("K47").Value = 25000
$ws.Range("M47").Value = -24028
$ws.Range("H96").Value = 2074.2666
$ws.Range("J96").Value = 2398.9167
$ws.Range("L96").Value = 7196.750100000001
$ws.Range("N96").Value = -9942.750100000001
$ws.Range("H107").Value = 425.18182
$ws.Range("I107").Value = 443.72223
$ws.Range("K107").Value = 443.72223
$ws.Range("M107").Value = 1476.27777
$ws.Range("H115").Value = 473.44446
$ws.Range("I115").Value = 473.44446
$ws.Range("K115").Value = 1420.33338
$ws.Range("M115").Value = 146.66662
$ws.Range("H137").Value = 6179.552
$ws.Range("I137").Value = 5600.3335
$ws.Range("K137").Value = 16801.0005
$ws.Range("M137").Value = -14251.0005
$ws.Range("H138").Value = 1892.7241
$ws.Range("J138").Value = 2059.1667
$ws.Range("L138").Value = 6177.500100000001
$ws.Range("N138").Value = -16457.5001
$ws.Range("H139").Value = 72000
$ws.Range("J139").Value = 72000
$ws.Range("L139").Value = 72000
$ws.Range("N139").Value = -82280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 73598
$ws.Range("J7").Value = 69330
$ws.Range("L7").Value = 69330
$ws.Range("N7").Value = -69558
$ws.Range("H45").Value = 23811896
$ws.Range("I45").Value = 35715996
$ws.Range("J45").Value = 3699
$ws.Range("K45").Value = 35715996
$ws.Range("L45").Value = 3699
$ws.Range("M45").Value = -35715619
$ws.Range("N45").Value = -4453
$ws.Range("H74").Value = 12382438
$ws.Range("I74").Value = 15626341
$ws.Range("K74").Value = 15626341
$ws.Range("M74").Value = -15625467
$ws.Range("H77").Value = 12382438
$ws.Range("I77").Value = 15626341
$ws.Range("K77").Value = 78131705
$ws.Range("M77").Value = -78127337

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3186.85
$ws.Range("I20").Value = 3785
$ws.Range("K20").Value = 3785
$ws.Range("M20").Value = -3538
$ws.Range("H96").Value = 46545.668
$ws.Range("I96").Value = 14832.333
$ws.Range("K96").Value = 14832.333
$ws.Range("M96").Value = -12086.333
$ws.Range("H99").Value = 2331.0513
$ws.Range("I99").Value = 1829.2667
$ws.Range("J99").Value = 4003.6667
$ws.Range("K99").Value = 1829.2667
$ws.Range("L99").Value = 4003.6667
$ws.Range("M99").Value = -331.2666999999999
$ws.Range("N99").Value = -6999.6667
$ws.Range("H117").Value = 99989.664
$ws.Range("J117").Value = 99989.664
$ws.Range("L117").Value = 99989.664
$ws.Range("N117").Value = -109167.664
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 109999.664
$ws.Range("J2").Value = 109999
$ws.Range("L2").Value = 109999
$ws.Range("N2").Value = -110225
$ws.Range("H99").Value = 3174.7144
$ws.Range("I99").Value = 2943.4
$ws.Range("K99").Value = 2943.4
$ws.Range("M99").Value = -1445.4
$ws.Range("H107").Value = 2642.9412
$ws.Range("I107").Value = 1340.2222
$ws.Range("K107").Value = 1340.2222
$ws.Range("M107").Value = 579.7778000000001
$ws.Range("H111").Value = 74989.5
$ws.Range("J111").Value = 74989.5
$ws.Range("L111").Value = 74989.5
$ws.Range("N111").Value = -83169.5
$ws.Range("H126").Value = 3174.7144
$ws.Range("I126").Value = 2943.4
$ws.Range("K126").Value = 8830.200000000001
$ws.Range("M126").Value = -6360.200000000001
$ws.Range("H134").Value = 3041.6428
$ws.Range("I134").Value = 2382.0417
$ws.Range("K134").Value = 7146.125100000001
$ws.Range("M134").Value = -4611.125100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 117.7619
$ws.Range("I2").Value = 79.0625
$ws.Range("J2").Value = 141.57692
$ws.Range("K2").Value = 474.375
$ws.Range("L2").Value = 849.4615200000001
$ws.Range("M2").Value = -361.375
$ws.Range("N2").Value = -1075.46152
$ws.Range("H8").Value = 187.9
$ws.Range("I8").Value = 187.9
$ws.Range("K8").Value = 563.7
$ws.Range("M8").Value = -424.7
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H122").Value = 3637.7144
$ws.Range("I122").Value = 790
$ws.Range("J122").Value = 6485.4287
$ws.Range("K122").Value = 7110
$ws.Range("L122").Value = 58368.85830000001
$ws.Range("M122").Value = -4660
$ws.Range("N122").Value = -63268.85830000001
$ws.Range("H131").Value = 296829.28
$ws.Range("I131").Value = 1668131.6
$ws.Range("J131").Value = 22568.8
$ws.Range("K131").Value = 5004394.800000001
$ws.Range("L131").Value = 67706.39999999999
$ws.Range("M131").Value = -4999354.800000001
$ws.Range("N131").Value = -77786.39999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3156
$ws.Range("I80").Value = 2853.0952
$ws.Range("J80").Value = 4064.7144
$ws.Range("K80").Value = 2853.0952
$ws.Range("L80").Value = 4064.7144
$ws.Range("M80").Value = -1855.0952
$ws.Range("N80").Value = -6060.7144
$ws.Range("H83").Value = 3156
$ws.Range("I83").Value = 2853.0952
$ws.Range("J83").Value = 4064.7144
$ws.Range("K83").Value = 14265.476
$ws.Range("L83").Value = 20323.572
$ws.Range("M83").Value = -9273.476000000001
$ws.Range("N83").Value = -30307.572
$ws.Range("H102").Value = 2324.606
$ws.Range("I102").Value = 1738.5834
$ws.Range("K102").Value = 1738.5834
$ws.Range("M102").Value = -116.5834
$ws.Range("H104").Value = 70833
$ws.Range("J104").Value = 70833
$ws.Range("L104").Value = 70833
$ws.Range("N104").Value = -77821
$ws.Range("H113").Value = 4579.3335
$ws.Range("I113").Value = 3100.5
$ws.Range("J113").Value = 5001.857
$ws.Range("K113").Value = 3100.5
$ws.Range("L113").Value = 5001.857
$ws.Range("M113").Value = -930.5
$ws.Range("N113").Value = -9341.857
$ws.Range("H117").Value = 105000
$ws.Range("J117").Value = 105000
$ws.Range("L117").Value = 105000
$ws.Range("N117").Value = -111884
$ws.Range("H122").Value = 1384.174
$ws.Range("I122").Value = 1421.5264
$ws.Range("K122").Value = 4264.5792
$ws.Range("M122").Value = -1814.5792
$ws.Range("H135").Value = 77890
$ws.Range("J135").Value = 77890
$ws.Range("L135").Value = 77890
$ws.Range("N135").Value = -88030

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H105").Value = 47333
$ws.Range("J105").Value = 47333
$ws.Range("L105").Value = 47333
$ws.Range("N105").Value = -54321
$ws.Range("H115").Value = 170000
$ws.Range("J115").Value = 170000
$ws.Range("L115").Value = 170000
$ws.Range("N115").Value = -172350
$ws.Range("H132").Value = 1102276
$ws.Range("I132").Value = 113639.445
$ws.Range("J132").Value = 10000005
$ws.Range("K132").Value = 340918.335
$ws.Range("L132").Value = 30000015
$ws.Range("M132").Value = -338388.335
$ws.Range("N132").Value = -30005075
$ws.Range("H136").Value = 142046.19
$ws.Range("I136").Value = 31875.75
$ws.Range("J136").Value = 205000.72
$ws.Range("K136").Value = 95627.25
$ws.Range("L136").Value = 615002.16
$ws.Range("M136").Value = -93077.25
$ws.Range("N136").Value = -620102.16
$ws.Range("H141").Value = 72571.664
$ws.Range("J141").Value = 72571.664
$ws.Range("L141").Value = 72571.664
$ws.Range("N141").Value = -82931.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2875.1516
$ws.Range("I122").Value = 2882.1738
$ws.Range("K122").Value = 8646.5214
$ws.Range("M122").Value = -6196.5214
$ws.Range("H125").Value = 98082
$ws.Range("J125").Value = 98082
$ws.Range("L125").Value = 98082
$ws.Range("N125").Value = -107922
